$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing existing row 2 (Ana) down to row 3
$ws.Rows.Item(2).Insert()
$ws.Range("A2:E2").ClearFormats()

$ws.Range("A2").Value = 101
$ws.Range("B2").Value = "Joao"
$ws.Range("C2").Value = "{101, 102, 103}"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35180184"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = 234.23
